# Update gh-pages to output generated at 456a3b4
# This script updates the "想去人数" (want-to-go count) column (F) values
# across the 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types)
# sheets of the 杭州-漫展信息 workbook to refresh them to newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet "展览") ---
$ws1.Range("F2").Value = 1217
$ws1.Range("F3").Value = 1114
$ws1.Range("F5").Value = 98
$ws1.Range("F6").Value = 58
$ws1.Range("F7").Value = 632
$ws1.Range("F8").Value = 79
$ws1.Range("F11").Value = 2202
$ws1.Range("F12").Value = 1547
$ws1.Range("F13").Value = 1229
$ws1.Range("F15").Value = 220
$ws1.Range("F16").Value = 482
$ws1.Range("F17").Value = 710
$ws1.Range("F18").Value = 262
$ws1.Range("F19").Value = 1079
$ws1.Range("F22").Value = 4146
$ws1.Range("F24").Value = 142
$ws1.Range("F25").Value = 117
$ws1.Range("F26").Value = 185
$ws1.Range("F28").Value = 596
$ws1.Range("F30").Value = 59
$ws1.Range("F33").Value = 353
$ws1.Range("F34").Value = 917
$ws1.Range("F35").Value = 117
$ws1.Range("F36").Value = 86
$ws1.Range("F37").Value = 109
$ws1.Range("F38").Value = 103

# --- 演出 (sheet "演出") ---
$ws2.Range("F3").Value = 773

# --- 全部类型 (sheet "全部类型") ---
$ws4.Range("F2").Value = 1217
$ws4.Range("F4").Value = 773
$ws4.Range("F5").Value = 1114
$ws4.Range("F9").Value = 98
$ws4.Range("F10").Value = 58
$ws4.Range("F11").Value = 632
$ws4.Range("F12").Value = 79
$ws4.Range("F16").Value = 2202
$ws4.Range("F17").Value = 1547
$ws4.Range("F18").Value = 1229
$ws4.Range("F20").Value = 220
$ws4.Range("F21").Value = 482
$ws4.Range("F23").Value = 710
$ws4.Range("F24").Value = 262
$ws4.Range("F25").Value = 1079
$ws4.Range("F28").Value = 4146
$ws4.Range("F30").Value = 142
$ws4.Range("F31").Value = 117
$ws4.Range("F32").Value = 185
$ws4.Range("F34").Value = 596
$ws4.Range("F36").Value = 59
$ws4.Range("F39").Value = 353
$ws4.Range("F40").Value = 917
$ws4.Range("F41").Value = 117
$ws4.Range("F42").Value = 86
$ws4.Range("F43").Value = 109
$ws4.Range("F44").Value = 103

$wb.Save()
